$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a 4-row x 20-column array holding the new table contents for rows 2..5.
$arr = New-Object 'object[,]' 4,20

# Row 2 (index 0): "Especial" quality, updated date/prices/origin -> Provincia de Quillota
$arr[0,0]  = 9
$arr[0,1]  = "Vega Central Mapocho de Santiago"
$arr[0,2]  = "Metropolitana"
$arr[0,3]  = 44915
$arr[0,4]  = 13
$arr[0,5]  = "Fruta"
$arr[0,6]  = 100104
$arr[0,7]  = "Frutos de pepita"
$arr[0,8]  = 100104004
$arr[0,9]  = "Níspero"
$arr[0,10] = "Golden Nugget"
$arr[0,11] = "Especial"
$arr[0,12] = 150
$arr[0,13] = 6000
$arr[0,14] = 6000
$arr[0,15] = 6000
$arr[0,16] = "$/bandeja 5 kilos"
$arr[0,17] = "Provincia de Quillota"
$arr[0,18] = 1200
$arr[0,19] = 5

# Row 3 (index 1): "Primera" quality, updated date/prices/origin -> Provincia de Quillota
$arr[1,0]  = 9
$arr[1,1]  = "Vega Central Mapocho de Santiago"
$arr[1,2]  = "Metropolitana"
$arr[1,3]  = 44915
$arr[1,4]  = 13
$arr[1,5]  = "Fruta"
$arr[1,6]  = 100104
$arr[1,7]  = "Frutos de pepita"
$arr[1,8]  = 100104004
$arr[1,9]  = "Níspero"
$arr[1,10] = "Golden Nugget"
$arr[1,11] = "Primera"
$arr[1,12] = 200
$arr[1,13] = 5000
$arr[1,14] = 5000
$arr[1,15] = 5000
$arr[1,16] = "$/bandeja 5 kilos"
$arr[1,17] = "Provincia de Quillota"
$arr[1,18] = 1000
$arr[1,19] = 5

# Row 4 (index 2): new row, same as original row 2 (Primera, O'Higgins)
$arr[2,0]  = 9
$arr[2,1]  = "Vega Central Mapocho de Santiago"
$arr[2,2]  = "Metropolitana"
$arr[2,3]  = 44911
$arr[2,4]  = 13
$arr[2,5]  = "Fruta"
$arr[2,6]  = 100104
$arr[2,7]  = "Frutos de pepita"
$arr[2,8]  = 100104004
$arr[2,9]  = "Níspero"
$arr[2,10] = "Golden Nugget"
$arr[2,11] = "Primera"
$arr[2,12] = 220
$arr[2,13] = 5000
$arr[2,14] = 5000
$arr[2,15] = 5000
$arr[2,16] = "$/bandeja 5 kilos"
$arr[2,17] = "Región de O'Higgins"
$arr[2,18] = 1000
$arr[2,19] = 5

# Row 5 (index 3): new row, same as original row 3 (Segunda, O'Higgins)
$arr[3,0]  = 9
$arr[3,1]  = "Vega Central Mapocho de Santiago"
$arr[3,2]  = "Metropolitana"
$arr[3,3]  = 44911
$arr[3,4]  = 13
$arr[3,5]  = "Fruta"
$arr[3,6]  = 100104
$arr[3,7]  = "Frutos de pepita"
$arr[3,8]  = 100104004
$arr[3,9]  = "Níspero"
$arr[3,10] = "Golden Nugget"
$arr[3,11] = "Segunda"
$arr[3,12] = 200
$arr[3,13] = 4000
$arr[3,14] = 4000
$arr[3,15] = 4000
$arr[3,16] = "$/bandeja 5 kilos"
$arr[3,17] = "Región de O'Higgins"
$arr[3,18] = 800
$arr[3,19] = 5

$ws.Range("A2:T5").Value = $arr

# Rows 4 and 5 are brand new cells; give column D the same date styling/number
# format used by the existing D2/D3 cells (style index "2", a date format).
$dateFormat = $ws.Range("D2").NumberFormat
$ws.Range("D2:D5").NumberFormat = $dateFormat

Write-Host "Updated range:" $ws.UsedRange.Address()
